# Updates cryptos list: prices and 1h volume % changes (and two name/row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.328.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").Value = "'2.961.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'520.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "'129.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'2.959.15"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").Value = "'0.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").Value = "'6.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("D11").Value = "'0.146"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.65%  "
$ws.Range("D12").Value = "'0.433"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("D14").Value = "'32.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "'3.434.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'60.239.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").Value = "'2.953.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "'6.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").Value = "'452.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.00%  "
$ws.Range("D21").Value = "'12.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "'0.661"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.70%  "
$ws.Range("D23").Value = "'6.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("D24").Value = "'77.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").Value = "'11.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "'7.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.51%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").Value = "'24.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("D32").Value = "'1.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("D33").Value = "'54.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'5.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D35").Value = "'2.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.34%  "
$ws.Range("D36").Value = "'5.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.87%  "
$ws.Range("D37").Value = "'445.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.30%  "
$ws.Range("D38").Value = "'3.146.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D39").Value = "'0.0767"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "'0.0373"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Value = "'7.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Value = "'2.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.38%  "
$ws.Range("D45").Value = "'0.240"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "'25.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.41%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.107"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'116.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").Value = "'1.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.68%  "
$ws.Range("D50").Value = "'0.0₃0495"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.41%  "
$ws.Range("D51").Value = "'1.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.07%  "
